$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "211.17", "15.30") are preserved exactly instead of being coerced
# into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.679.26"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.598.72"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "211.17"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "19.71"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.822.01"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "1.600.04"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "64.86"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "26.654.19"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "209.92"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "146.32"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").Value = "15.30"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "1.296.60"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "0.842"
$ws.Range("E39").Value = "  +2.86%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "63.97"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "1.734.94"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "0.895"
$ws.Range("E46").Value = "  +10.94%  "
$ws.Range("D47").Value = "90.06"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").Value = "7.48"
$ws.Range("E51").Value = "  +0.90%  "

# Restore original (default) formatting on column D now that the
# text values have been written, so cell styles match the source file.
$ws.Range("D2:D51").ClearFormats()
